$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3, duplicating row 2's data with a different name
# (written before row 2's rename so shared-string order matches)
$ws.Range("A3").Value2 = "Dawn of War 1"

# Update existing row 2: rename product and mark status active
$ws.Range("A2").Value2 = "Dawn of War 2"
$ws.Range("I2").Value2 = 1

$ws.Range("B3").Value2 = $ws.Range("B2").Value2
$ws.Range("C3").Value2 = $ws.Range("C2").Value2
$ws.Range("D3").Value2 = $ws.Range("D2").Value2
$ws.Range("E3").Value2 = $ws.Range("E2").Value2
$ws.Range("F3").Value2 = $ws.Range("F2").Value2
$ws.Range("G3").Value2 = $ws.Range("G2").Value2
$ws.Range("H3").Value2 = $ws.Range("H2").Value2
$ws.Range("I3").Value2 = 1
$ws.Range("J3").Value2 = $false
$ws.Range("K3").Value2 = $false

$ws.Range("A3").Select() | Out-Null
